$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Rename "Sort kaylee photo" -> "Sort Kathryn photo" (H10), and add new task
# "transfer jpeg to imgur" (K10) for the About Us row group.
$ws1.Range("H10").Value = "Sort Kathryn photo"
$ws1.Range("K10").Value = "transfer jpeg to imgur"

# Fix typo / trailing "?" in the google maps task text.
$ws1.Range("H23").Value = "add in actual google maps api"

# New task under Switzerland: "add marker to google maps page"
$ws1.Range("H24").Value = "add marker to google maps page"

# Move the "Overall" section (F25/H25) down to rows 27 (F27/H27)
$ws1.Range("F25").Value = $null
$ws1.Range("H25").Value = $null
$ws1.Range("F27").Value = "Overall"
$ws1.Range("H27").Value = "Sort CSS"

# New "Navigation" section added at rows 29-30
$ws1.Range("F29").Value = "Navigation"
$ws1.Range("H29").Value = "edit on destinations pages to look more suited"
$ws1.Range("H30").Value = "dropdown menu on destinations menu bit?"

# Apply green fill highlight style to the "done/completed" style cells
$doneCells = @("H9", "H17", "H18", "H19", "H22", "H23", "H24")
foreach ($addr in $doneCells) {
    $ws1.Range($addr).Interior.Color = 5287936
}

# Update selection to match the saved state
$ws1.Range("J27").Select()
